$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Septiembre de 2020 a las 12:54"

# --- Row 4 ---
$ws.Range("B4").Value = 7185915
$ws.Range("C4").Value = 444
$ws.Range("D4").Value = 4438628
$ws.Range("E4").Value = 2539747
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 207540

# --- Row 5 ---
$ws.Range("B5").Value = 5823060
$ws.Range("C5").Value = 6957
$ws.Range("E5").Value = 974549
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 92347

# --- Row 16 ---
$ws.Range("B16").Value = 439882
$ws.Range("C16").Value = 3563
$ws.Range("D16").Value = 369842
$ws.Range("E16").Value = 44818
$ws.Range("G16").Value = 207
$ws.Range("H16").Value = 25222

# --- Row 32 ---
$ws.Range("B32").Value = 124650
$ws.Range("C32").Value = 225
$ws.Range("D32").Value = 121512
$ws.Range("E32").Value = 2926

# --- Row 33 ---
$ws.Range("B33").Value = 119683
$ws.Range("C33").Value = 1629
$ws.Range("D33").Value = 96158
$ws.Range("E33").Value = 18892
$ws.Range("G33").Value = 42
$ws.Range("H33").Value = 4633

# --- Row 44 ---
$ws.Range("B44").Value = 89540
$ws.Range("C44").Value = 1008
$ws.Range("D44").Value = 78819
$ws.Range("E44").Value = 10312
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 409

# --- Row 49 ---
$ws.Range("B49").Value = 76957
$ws.Range("C49").Value = 306
$ws.Range("D49").Value = 73925
$ws.Range("E49").Value = 2225
$ws.Range("G49").Value = 5
$ws.Range("H49").Value = 807

# --- Row 56 ---
$ws.Range("E56").Value = 6614
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 234

# --- Row 61 ---
$ws.Range("B61").Value = 51864
$ws.Range("C61").Value = 372
$ws.Range("E61").Value = 7201
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 2063

# --- Row 79 ---
$ws.Range("B79").Value = 26564
$ws.Range("C79").Value = 248
$ws.Range("D79").Value = 19507
$ws.Range("E79").Value = 6249
$ws.Range("G79").Value = 6
$ws.Range("H79").Value = 808

# --- Row 97 ---
$ws.Range("B97").Value = 10687
$ws.Range("C97").Value = 111
$ws.Range("D97").Value = 9696
$ws.Range("E97").Value = 858

# --- Row 109 ---
$ws.Range("A109").Value = "Eslovaquia"
$ws.Range("B109").Value = 8048
$ws.Range("C109").Value = 419
$ws.Range("D109").Value = 4036
$ws.Range("E109").Value = 3971
$ws.Range("H109").Value = 41

# --- Row 110 ---
$ws.Range("A110").Value = "Zimbabue"
$ws.Range("B110").Value = 7752
$ws.Range("D110").Value = 6043
$ws.Range("E110").Value = 1482
$ws.Range("H110").Value = 227

# --- Row 142 ---
$ws.Range("B142").Value = 3334
$ws.Range("C142").Value = 1
$ws.Range("D142").Value = 3158
$ws.Range("E142").Value = 163

# --- Row 145 ---
$ws.Range("A145").Value = "Malta"
$ws.Range("B145").Value = 2929
$ws.Range("C145").Value = 31
$ws.Range("D145").Value = 2261
$ws.Range("E145").Value = 639
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 29

# --- Row 146 ---
$ws.Range("A146").Value = "Botsuana"
$ws.Range("B146").Value = 2921
$ws.Range("D146").Value = 701
$ws.Range("E146").Value = 2204
$ws.Range("H146").Value = 16

# --- Row 215 ---
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

# --- Row 216 ---
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
